# Actualizacion lista de precios mar 16/09/2025 22:01:28,76
#
# The price list ("Lista") had the price in C2 updated from 3 to 100,
# and the active cell selection was left on C2 (it had been on D4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the price in C2.
$ws.Range("C2").Value = 100

# Leave the selection on C2.
$ws.Range("C2").Select()
